$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the daily price log for
# "Vega Modelo de Temuco - Limón". The new record belongs right after
# row 881 (chronologically among the existing rows), so insert a new
# row at position 882 which pushes the former rows 882-938 down to
# become rows 883-939, and populate the newly inserted row with the
# new reading.

$ws.Rows(882).Insert()

$ws.Range("A882").Value = 10
$ws.Range("B882").Value = "Vega Modelo de Temuco"
$ws.Range("C882").Value = "La Araucanía"
$ws.Range("D882").Value = 44516
$ws.Range("E882").Value = 9
$ws.Range("F882").Value = "Fruta"
$ws.Range("G882").Value = 100102
$ws.Range("H882").Value = "Cítricos"
$ws.Range("I882").Value = 100102003
$ws.Range("J882").Value = "Limón"
$ws.Range("K882").Value = "Sin especificar"
$ws.Range("L882").Value = "1a amarillo"
$ws.Range("M882").Value = 180
$ws.Range("N882").Value = 9000
$ws.Range("O882").Value = 10000
$ws.Range("P882").Value = 9361
$ws.Range("Q882").Value = "`$/bandeja 15 kilos"
$ws.Range("R882").Value = "Región de O'Higgins"
$ws.Range("S882").Value = 624
$ws.Range("T882").Value = 15
